$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 128; this shifts existing rows 128-140 down to 129-141
$ws.Rows.Item(128).Insert()

# Populate the newly inserted row 128 with the new data record
$ws.Range("A128").Value = 4
$ws.Range("B128").Value = 'Feria Lagunitas de Puerto Montt'
$ws.Range("C128").Value = 'Los Lagos'
$ws.Range("D128").Value = '9/10/2021'
$ws.Range("E128").Value = 10
$ws.Range("F128").Value = 100112021
$ws.Range("G128").Value = 'Ají'
$ws.Range("H128").Value = 'Inferno'
$ws.Range("I128").Value = 'Extra'
$ws.Range("J128").Value = 60
$ws.Range("K128").Value = 50000
$ws.Range("L128").Value = 50000
$ws.Range("M128").Value = 50000
$ws.Range("N128").Value = '$/caja 12 kilos'
$ws.Range("O128").Value = 'Región de Arica y Parinacota'
$ws.Range("P128").Value = 4167
$ws.Range("Q128").Value = 12
$ws.Range("R128").Value = 'Hortaliza'
